$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 275, pushing existing rows 275-331 down to 277-333
$ws.Rows("275:276").Insert()

# --- New row 275: Artic Mist / Primera ---
$ws.Range("A275").Value = 11
$ws.Range("B275").Value = "Vega Monumental Concepción"
$ws.Range("C275").Value = "Bíobío"
$ws.Range("D275").Value = 44642
$ws.Range("E275").Value = 8
$ws.Range("F275").Value = "Fruta"
$ws.Range("G275").Value = 100103
$ws.Range("H275").Value = "Frutos de hueso (carozo)"
$ws.Range("I275").Value = 100103006
$ws.Range("J275").Value = "Nectarín"
$ws.Range("K275").Value = "Artic Mist"
$ws.Range("L275").Value = "Primera"
$ws.Range("M275").Value = 220
$ws.Range("N275").Value = 12000
$ws.Range("O275").Value = 13000
$ws.Range("P275").Value = 12455
$ws.Range("Q275").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R275").Value = "Región de O'Higgins"
$ws.Range("S275").Value = 778
$ws.Range("T275").Value = 16

# --- New row 276: Ruby Diamond / Primera ---
$ws.Range("A276").Value = 11
$ws.Range("B276").Value = "Vega Monumental Concepción"
$ws.Range("C276").Value = "Bíobío"
$ws.Range("D276").Value = 44642
$ws.Range("E276").Value = 8
$ws.Range("F276").Value = "Fruta"
$ws.Range("G276").Value = 100103
$ws.Range("H276").Value = "Frutos de hueso (carozo)"
$ws.Range("I276").Value = 100103006
$ws.Range("J276").Value = "Nectarín"
$ws.Range("K276").Value = "Ruby Diamond"
$ws.Range("L276").Value = "Primera"
$ws.Range("M276").Value = 220
$ws.Range("N276").Value = 13000
$ws.Range("O276").Value = 14000
$ws.Range("P276").Value = 13455
$ws.Range("Q276").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R276").Value = "Región de O'Higgins"
$ws.Range("S276").Value = 841
$ws.Range("T276").Value = 16
